# technical_report.xlsx edit: "criacao do delete por id e listar todos"
#
# 1. Replace the sample report row (row 2) text with new content.
# 2. Make the header row (row 1) bold, 14pt Calibri.
# 3. Auto-fit the columns so the widened/bolded header still fits.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the data row -------------------------------------------------
$ws.Range("A2").Value = "Angelina"
$ws.Range("B2").Value = "Celular Android"
$ws.Range("C2").Value = "não liga"
$ws.Range("D2").Value = "caiu na agua"

# --- Bold + enlarge the header row ---------------------------------------
# Set both font attributes on a single cell first (keeps the style table
# tight: one new font / one new cell style), then copy that formatting
# across the rest of the header instead of re-deriving it cell by cell.
$headerFirst = $ws.Range("A1")
$headerFirst.Font.FontStyle = "Bold"
$headerFirst.Font.Size = 14

$headerFirst.Copy()
$ws.Range("B1:D1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Resize columns to fit the new header formatting ----------------------
$ws.Columns("A:D").AutoFit()
